$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new coordinates + formula for AD3 (distance calc now applies here too) ---
$ws.Range("AA3").Value = 446820
$ws.Range("AB3").Value = 7395454
$ws.Range("AD3").Formula = '=SQRT((AA3-$AA$7)^2+(AB3-$AB$7)^2)'

# --- Row 4: Record_start date unknown -> "-" ---
$ws.Range("G4").Value = "-"

# --- Row 8: fix the Record_start typo'd date (2001.20.03 -> 2001.03.20) ---
# Write it through a formula + paste-special-values so Excel stores it as plain
# text (matching the original authoring) instead of auto-converting the
# date-looking string into a date serial number.
$r = $ws.Range("G8")
$r.Formula = '="2001.03.20"'
$r.Copy()
$r.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Sheet view: split pane + selections ---
$ws.Range("A10:XFD10").Select()
$ws.Range("A10").Activate()
$win = $excel.ActiveWindow
$win.SplitColumn = 2
$win.Split = $true
$win.Panes.Item(2).Activate()
$ws.Range("G9").Select()

# --- Workbook window size ---
$win.Height = 17060
